# Add MPs' base pay (Freddo) multiplier row and CPIH multiplier row to the
# "Inflation" sheet, shifting the existing CPI/RPI multiplier rows down by
# one, and extend the trailing blank spacer rows by one row to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inflation")

# --- Row 9 (new): CPIH Multiplier, derived from row 5 (CPIH index) -------
# Written first while row 9 is still blank/unused.
$ws.Range("A9").Value2 = "CPIH Multiplier"
$ws.Range("B9").Formula = '=$B$5/B5'
$ws.Range("C9:Q9").Formula = '=$B$5/C5'

# --- Row 8: RPI Multiplier, derived from row 4 (RPI index) ---------------
# This re-creates what used to live in row 7, now one row further down.
$ws.Range("A8").Value2 = "RPI Multiplier"
$ws.Range("A8").Font.Bold = $true
$ws.Range("B8").Formula = '=$B$4/B4'
$ws.Range("C8:Q8").Formula = '=$B$4/C4'

# --- Row 7: CPI Multiplier, derived from row 3 (CPI index) ---------------
# This re-creates what used to live in row 6, now one row further down.
$ws.Range("A7").Value2 = "CPI Multiplier"
$ws.Range("A7").Font.Bold = $true
$ws.Range("B7").Formula = '=$B$3/B3'
$ws.Range("C7:Q7").Formula = '=$B$3/C3'

# --- Row 6 (new): Freddo Multiplier, derived from row 2 (Freddo index) ---
$ws.Range("A6").Value2 = "Freddo Multiplier"
$ws.Range("A6").Font.Bold = $true
$ws.Range("B6").Formula = '=$B$2/B2'
$ws.Range("C6:Q6").Formula = '=$B$2/C2'

# --- Extend the trailing blank spacer rows by one (row 29) ----------------
$ws.Range("D29").Value2 = $ws.Range("D28").Value2
$ws.Range("D29").Style = $ws.Range("D28").Style

# --- Selection moves to G27 (was D27) -------------------------------------
$ws.Range("G27").Select()
